$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scenario row: code first (so it becomes shared string 137), then
# the description (so it becomes shared string 138) - matching the
# order the strings were appended to the shared-strings table.
$code = "statements:`n    - context: workMode := user[""urn:ietf:params:scim:schemas:extension:ibm:2.0:User""].customAttributes[0].values[0]`n    - return: context.workMode.toUpper()"
$description = "get a custom attribute called ""workMode"" in uppercase"

$ws.Range("B70").Value = $code
$ws.Range("A70").Value = $description

# Match the wrap-text style (cellXfs index 1) used by the other
# multi-line rows in the sheet.
$ws.Range("A70:B70").WrapText = $true

# Size the row to fit the 3-line code snippet (3 x 14.4pt line height).
$ws.Rows.Item(70).RowHeight = 43.2

# Move the active selection the way the author's workbook ended up.
$ws.Range("A72").Select()
